# Program Management: roll the "Data" sheet over to the new cartridge batch
# "TestCartridge0237" (sample matrix SMatrix_Salm_0237 / lab sample id
# 20220512-TestAut-PA-30237, result date 05/12/2022), replacing the values
# left over from the previous "TestCartridge5523" batch in rows 2-13, and
# renumbering the Lane column sequentially 1-12.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data")

$resultIds = @(
    "A1426201",
    "A1426202",
    "A1426203",
    "A1426204",
    "A1426205",
    "A1426206",
    "A1426207",
    "A1426208",
    "A1426209",
    "A1426210",
    "A1426211",
    "A1426212"
)

$sampleMatrix = "SMatrix_Salm_0237"
$labSampleId  = "20220512-TestAut-PA-30237"
$resultDate   = "05/12/2022"
$cartridgeId  = "TestCartridge0237"

for ($i = 0; $i -lt $resultIds.Length; $i++) {
    $row = $i + 2

    $ws.Cells.Item($row, 1).Value  = $resultIds[$i]      # A: Result ID
    $ws.Cells.Item($row, 4).Value  = $sampleMatrix       # D: Sample Matrix
    $ws.Cells.Item($row, 5).Value  = $labSampleId        # E: Lab Sample ID

    # Q (Lane) and R (Result Date) hold plain text that merely looks like a
    # number/date ("1".."12", "05/12/2022"). Force text format before the
    # assignment so Excel doesn't auto-convert them to a number/date
    # serial, then drop back to the original (General) formatting.
    $lane = $ws.Cells.Item($row, 17)
    $lane.NumberFormat = "@"
    $lane.Value = [string]($i + 1)                       # Q: Lane (1..12)
    $lane.ClearFormats()

    $result_date = $ws.Cells.Item($row, 18)
    $result_date.NumberFormat = "@"
    $result_date.Value = $resultDate                     # R: Result Date
    $result_date.ClearFormats()

    $ws.Cells.Item($row, 20).Value = $cartridgeId        # T: Cartridge ID
}
